$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.581.82'
$ws.Range("E2").Value = '  -0.73%  '
$ws.Range("D3").Value = '1.860.95'
$ws.Range("E3").Value = '  -1.01%  '
$ws.Range("D4").Value = '1.012'
$ws.Range("E4").Value = '  +0.83%  '
$ws.Range("D5").Value = '333.37'
$ws.Range("E5").Value = '  +0.16%  '
$ws.Range("D6").Value = '1.011'
$ws.Range("E6").Value = '  +0.77%  '
$ws.Range("D7").Value = '0.4663'
$ws.Range("E7").Value = '  -1.27%  '
$ws.Range("D8").Value = '0.3886'
$ws.Range("E8").Value = '  -1.78%  '
$ws.Range("D9").Value = '45.80'
$ws.Range("E9").Value = '  -4.39%  '
$ws.Range("E10").Value = '  -1.31%  '
$ws.Range("D11").Value = '1.000'
$ws.Range("E11").Value = '  -3.16%  '
$ws.Range("D12").Value = '21.66'
$ws.Range("E12").Value = '  -2.82%  '
$ws.Range("D13").Value = '1.866.45'
$ws.Range("E13").Value = '  -0.70%  '
$ws.Range("D14").Value = '5.977'
$ws.Range("D15").Value = '7.199'
$ws.Range("E15").Value = '  +0.92%  '
$ws.Range("E16").Value = '  +0.68%  '
$ws.Range("D17").Value = '87.77'
$ws.Range("E17").Value = '  +0.55%  '
$ws.Range("D18").Value = '0.06695'
$ws.Range("E18").Value = '  +0.45%  '
$ws.Range("D19").Value = '0.00001043'
$ws.Range("E19").Value = '  -0.32%  '
$ws.Range("D20").Value = '16.89'
$ws.Range("E20").Value = '  -2.45%  '
$ws.Range("D21").Value = '1.011'
$ws.Range("E21").Value = '  +0.79%  '
$ws.Range("D22").Value = '27.558.06'
$ws.Range("E22").Value = '  -0.86%  '
$ws.Range("D23").Value = '5.434'
$ws.Range("E23").Value = '  -1.86%  '
$ws.Range("E24").Value = '  -1.77%  '
$ws.Range("D25").Value = '2.307'
$ws.Range("E25").Value = '  -0.01%  '
$ws.Range("D26").Value = '2.086.89'
$ws.Range("E26").Value = '  -0.61%  '
$ws.Range("D27").Value = '158.98'
$ws.Range("E27").Value = '  -0.32%  '
$ws.Range("D28").Value = '19.70'
$ws.Range("E28").Value = '  -2.50%  '
$ws.Range("D29").Value = '2.113'
$ws.Range("E29").Value = '  +0.25%  '
$ws.Range("D30").Value = '5.379'
$ws.Range("E30").Value = '  -3.67%  '
$ws.Range("D31").Value = '121.13'
$ws.Range("E31").Value = '  -0.67%  '
$ws.Range("D32").Value = '0.9719'
$ws.Range("E32").Value = '  -1.11%  '
$ws.Range("D33").Value = '0.09458'
$ws.Range("E33").Value = '  -0.88%  '
$ws.Range("D34").Value = '3.651'
$ws.Range("E34").Value = '  +1.55%  '
$ws.Range("D35").Value = '5.287'
$ws.Range("E35").Value = '  -1.48%  '
$ws.Range("D36").Value = '1.329'
$ws.Range("E36").Value = '  -8.35%  '
$ws.Range("D37").Value = '0.06021'
$ws.Range("E37").Value = '  -1.81%  '
$ws.Range("D38").Value = '0.02210'
$ws.Range("E38").Value = '  -2.13%  '
$ws.Range("D39").Value = '1.193'
$ws.Range("E39").Value = '  -3.44%  '
$ws.Range("D40").Value = '8.168'
$ws.Range("E40").Value = '  +0.36%  '
$ws.Range("D42").Value = '0.5901'
$ws.Range("E42").Value = '  -2.21%  '
$ws.Range("D43").Value = '0.1876'
$ws.Range("E43").Value = '  -1.52%  '
$ws.Range("D44").Value = '10.21'
$ws.Range("E44").Value = '  -0.63%  '
$ws.Range("D45").Value = '1.241'
$ws.Range("E45").Value = '  -2.69%  '
$ws.Range("D46").Value = '0.5609'
$ws.Range("E46").Value = '  -2.32%  '
$ws.Range("D47").Value = '12.12'
$ws.Range("E47").Value = '  -0.92%  '
$ws.Range("E48").Value = '  -2.15%  '
$ws.Range("D49").Value = '3.274'
$ws.Range("E49").Value = '  -2.99%  '
$ws.Range("D50").Value = '0.06759'
$ws.Range("E50").Value = '  -2.30%  '
$ws.Range("D51").Value = '112.36'
$ws.Range("E51").Value = '  -1.80%  '
